# Rename "Investor" header to "Stakeholder" and "PAN" header to "PAN/Tax Id",
# then reset the view (top-left cell / selection) back to A1.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update header row text
$ws.Range("A1").Value = "Stakeholder"
$ws.Range("D1").Value = "PAN/Tax Id"

# Reset the scroll position / selection back to A1 (was topLeftCell D1, selection V6)
$excel.ActiveWindow.ScrollRow = 1
$excel.ActiveWindow.ScrollColumn = 1
$ws.Range("A1").Select()
